$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.192.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.26%  '

$ws.Range("D3").Value = '''1.681.74'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.32%  '

$ws.Range("D4").Value = '''1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''216.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.52%  '

$ws.Range("D6").Value = '''0.5246'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.77%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '''0.2691'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.52%  '

$ws.Range("D9").Value = '''0.06364'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.65%  '

$ws.Range("E10").Value = '  -1.84%  '

$ws.Range("D11").Value = '''0.07636'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.51%  '

$ws.Range("D12").Value = '''1.693.19'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.94%  '

$ws.Range("D13").Value = '''4.517'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.13%  '

$ws.Range("D14").Value = '''0.5749'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.26%  '

$ws.Range("D15").Value = '''0.000008302'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.92%  '

$ws.Range("D16").Value = '''65.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.20%  '

$ws.Range("D17").Value = '''26.251.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("D18").Value = '''1.006'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("D19").Value = '''4.857'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.10%  '

$ws.Range("D20").Value = '''10.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.47%  '

$ws.Range("D21").Value = '''188.91'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.62%  '

$ws.Range("D22").Value = '''6.238'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.84%  '

$ws.Range("D23").Value = '''1.007'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").Value = '''148.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.47%  '

$ws.Range("D25").Value = '''7.782'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.36%  '

$ws.Range("D26").Value = '''0.1260'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.38%  '

$ws.Range("E27").Value = '  -0.09%  '

$ws.Range("D28").Value = '''0.06312'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.28%  '

$ws.Range("D29").Value = '''1.377'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.00%  '

$ws.Range("D30").Value = '''1.315'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.25%  '

$ws.Range("D31").Value = '''3.565'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.44%  '

$ws.Range("D32").Value = '''3.567'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.50%  '

$ws.Range("E33").Value = '  +1.50%  '

$ws.Range("E34").Value = '  -0.51%  '

$ws.Range("D35").Value = '''0.6112'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.14%  '

$ws.Range("E36").Value = '  +0.61%  '

$ws.Range("D37").Value = '''2.755'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.23%  '

$ws.Range("D38").Value = '''6.193'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.10%  '

$ws.Range("D39").Value = '''0.01616'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.35%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''0.8912'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.99%  '

$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '''1.094.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.79%  '

$ws.Range("D42").Value = '''1.010'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.42%  '

$ws.Range("D43").Value = '''100.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("D44").Value = '''1.832.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.34%  '

$ws.Range("E45").Value = '  -0.73%  '

$ws.Range("D46").Value = '''57.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.59%  '

$ws.Range("D47").Value = '''1.007'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.71%  '

$ws.Range("D48").Value = '''8.058'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.54%  '

$ws.Range("D49").Value = '''0.05279'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.35%  '

$ws.Range("D50").Value = '''0.4282'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.15%  '

$ws.Range("D51").Value = '''6.012'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.08%  '
